$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '69.340.17'
$ws.Range("E2").Value = '  -0.03%  '

# Row 3
$ws.Range("D3").Value = '3.689.97'
$ws.Range("E3").Value = '  +0.12%  '

# Row 4
$ws.Range("E4").Value = '  +0.06%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '680.77'
$ws.Range("E5").Value = '  -1.10%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '159.24'
$ws.Range("E6").Value = '  -1.57%  '

# Row 7
$ws.Range("E7").Value = '  -0.03%  '

# Row 8
$ws.Range("E8").Value = '  -0.42%  '

# Row 9
$ws.Range("E9").Value = '  -0.85%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.12'
$ws.Range("E10").Value = '  -3.44%  '

# Row 11
$ws.Range("E11").Value = '  -0.04%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000231'
$ws.Range("E12").Value = '  -2.32%  '

# Row 13
$ws.Range("D13").Value = '4.313.83'
$ws.Range("E13").Value = '  +0.16%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.43'
$ws.Range("E14").Value = '  -1.99%  '

# Row 15
$ws.Range("D15").Value = '3.685.40'
$ws.Range("E15").Value = '  +0.01%  '

# Row 16
$ws.Range("D16").Value = '69.331.20'

# Row 17
$ws.Range("E17").Value = '  +2.07%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '16.01'
$ws.Range("E18").Value = '  -0.41%  '

# Row 19
$ws.Range("E19").Value = '  -0.56%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '468.75'
$ws.Range("E20").Value = '  -1.63%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.91'
$ws.Range("E21").Value = '  -0.36%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.654'
$ws.Range("E22").Value = '  -0.72%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '80.00'
$ws.Range("E23").Value = '  +0.20%  '

# Row 24
$ws.Range("D24").Value = '3.837.41'
$ws.Range("E24").Value = '  +0.12%  '

# Row 26
$ws.Range("E26").Value = '  -4.68%  '

# Row 27
$ws.Range("E27").Value = '  -3.10%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.14'
$ws.Range("E28").Value = '  -2.53%  '

# Row 29
$ws.Range("E29").Value = '  -0.58%  '

# Row 30
$ws.Range("E30").Value = '  -3.72%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.58'
$ws.Range("E31").Value = '  -3.07%  '

# Row 32
$ws.Range("E32").Value = '  -2.70%  '

# Row 33
$ws.Range("E33").Value = '  +0.17%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '26.95'
$ws.Range("E34").Value = '  +0.38%  '

# Row 35
$ws.Range("D35").Value = '3.678.94'
$ws.Range("E35").Value = '  +0.53%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.157'
$ws.Range("E36").Value = '  -6.02%  '

# Row 37
$ws.Range("E37").Value = '  -0.71%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.24'
$ws.Range("E38").Value = '  -0.64%  '

# Row 39
$ws.Range("E39").Value = '  -0.01%  '

# Row 40
$ws.Range("E40").Value = '  -0.05%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.21'
$ws.Range("E41").Value = '  -4.69%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0904'
$ws.Range("E42").Value = '  -1.65%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '171.35'
$ws.Range("E43").Value = '  +4.28%  '

# Row 44
$ws.Range("E44").Value = '  -0.87%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '47.53'
$ws.Range("E45").Value = '  -1.12%  '

# Row 46
$ws.Range("B46").Value = 'SuiNetwork'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.12'
$ws.Range("E46").Value = '  -0.93%  '

# Row 47
$ws.Range("B47").Value = 'InjectiveProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '28.25'
$ws.Range("E47").Value = '  -5.16%  '

# Row 48
$ws.Range("B48").Value = 'dogwifhat'
$ws.Range("C48").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.70'
$ws.Range("E48").Value = '  -1.42%  '

# Row 49
$ws.Range("B49").Value = 'ONDO'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.31'
$ws.Range("E49").Value = '  -0.82%  '

# Row 50
$ws.Range("E50").Value = '  -1.73%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.80'
$ws.Range("E51").Value = '  -2.82%  '
